$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for updated crypto data.
# D-column values that parse as plain numbers are written with a leading
# apostrophe so Excel keeps them as text (matching the source data's text type),
# exactly like the original "double-dot" price strings that cannot be parsed as numbers.

$ws.Range("D2").Value = "27.287.83"
$ws.Range("E2").Value = "  -2.56%  "
$ws.Range("D3").Value = "1.706.57"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'223.16"
$ws.Range("E5").Value = "  -2.80%  "
$ws.Range("D6").Value = "'0.5303"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.2654"
$ws.Range("E8").Value = "  -4.49%  "
$ws.Range("D9").Value = "'0.06583"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("D10").Value = "'20.87"
$ws.Range("E10").Value = "  -4.26%  "
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").Value = "'4.587"
$ws.Range("E12").Value = "  -2.60%  "
$ws.Range("D13").Value = "1.725.46"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").Value = "1.941.17"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "'0.5736"
$ws.Range("E15").Value = "  -4.51%  "
$ws.Range("D16").Value = "0.0₅8185"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "'67.48"
$ws.Range("E17").Value = "  -3.10%  "
$ws.Range("D18").Value = "27.268.28"
$ws.Range("E18").Value = "  -2.52%  "
$ws.Range("D19").Value = "'216.37"
$ws.Range("E19").Value = "  -3.27%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'4.677"
$ws.Range("E21").Value = "  -3.10%  "
$ws.Range("D22").Value = "'10.44"
$ws.Range("E22").Value = "  -4.61%  "
$ws.Range("D23").Value = "'5.975"
$ws.Range("E23").Value = "  -4.57%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "'142.23"
$ws.Range("E25").Value = "  -3.29%  "
$ws.Range("D26").Value = "'1.743"
$ws.Range("E26").Value = "  +5.60%  "
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("D28").Value = "'7.255"
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("D29").Value = "'16.32"
$ws.Range("E29").Value = "  -3.90%  "
$ws.Range("D30").Value = "'0.05371"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("D31").Value = "'1.291"
$ws.Range("E31").Value = "  -2.22%  "
$ws.Range("D32").Value = "'3.508"
$ws.Range("D33").Value = "'3.423"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").Value = "'1.637"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").Value = "'2.876"
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("D36").Value = "'2.419"
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("D37").Value = "'0.9467"
$ws.Range("E37").Value = "  -3.89%  "
$ws.Range("D38").Value = "'0.5883"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("D39").Value = "'0.01635"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("D40").Value = "'5.858"
$ws.Range("E40").Value = "  -1.24%  "
$ws.Range("D41").Value = "'1.002"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "1.039.79"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").Value = "'0.8399"
$ws.Range("E43").Value = "  -0.65%  "
$ws.Range("D44").Value = "'101.02"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").Value = "1.848.11"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("D47").Value = "'58.11"
$ws.Range("E47").Value = "  -3.26%  "
$ws.Range("D48").Value = "'0.4497"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").Value = "'1.008"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.086"
$ws.Range("E50").Value = "  -2.69%  "
$ws.Range("B51").Value = "XinFinNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D51").Value = "'0.06565"
$ws.Range("E51").Value = "  +10.45%  "
